$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 109.8125
$ws.Range("I6").Value = 109.8125
$ws.Range("K6").Value = 329.4375
$ws.Range("M6").Value = -217.4375

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""

$ws.Range("H32").Value = 2216.25
$ws.Range("J32").Value = 2370.7144
$ws.Range("L32").Value = 2370.7144
$ws.Range("N32").Value = -3022.7144

$ws.Range("H40").Value = 2400
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = ""

$ws.Range("H43").Value = 1249.75
$ws.Range("I43").Value = 999.6667
$ws.Range("K43").Value = 999.6667
$ws.Range("M43").Value = -930.6667

$ws.Range("H87").Value = 42000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 42000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 42000
$ws.Range("M87").Value = ""
$ws.Range("N87").Value = -44496

$ws.Range("H90").Value = 42000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 42000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 126000
$ws.Range("M90").Value = ""
$ws.Range("N90").Value = -138480

$ws.Range("H98").Value = 1612.5
$ws.Range("I98").Value = 1750
$ws.Range("K98").Value = 1750
$ws.Range("M98").Value = -252

$ws.Range("H122").Value = 1612.5
$ws.Range("I122").Value = 1750
$ws.Range("K122").Value = 5250
$ws.Range("M122").Value = -2800

$ws.Range("H137").Value = 2108.1177
$ws.Range("I137").Value = 1774.2858
$ws.Range("J137").Value = 3666
$ws.Range("K137").Value = 5322.857400000001
$ws.Range("L137").Value = 10998
$ws.Range("M137").Value = -2772.857400000001
$ws.Range("N137").Value = -16098

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1549.7333
$ws.Range("I31").Value = 1075.2
$ws.Range("J31").Value = 2498.8
$ws.Range("K31").Value = 1075.2
$ws.Range("L31").Value = 2498.8
$ws.Range("M31").Value = -780.2
$ws.Range("N31").Value = -3088.8

$ws.Range("H34").Value = 1549.7333
$ws.Range("I34").Value = 1075.2
$ws.Range("J34").Value = 2498.8
$ws.Range("K34").Value = 1075.2
$ws.Range("L34").Value = 2498.8
$ws.Range("M34").Value = -873.2
$ws.Range("N34").Value = -2902.8

$ws.Range("H58").Value = 4295.316
$ws.Range("I58").Value = 2704
$ws.Range("J58").Value = 4593.6875
$ws.Range("K58").Value = 2704
$ws.Range("L58").Value = 4593.6875
$ws.Range("M58").Value = -2501
$ws.Range("N58").Value = -4999.6875

$ws.Range("H59").Value = 90000
$ws.Range("J59").Value = 90000
$ws.Range("L59").Value = 90000
$ws.Range("N59").Value = -92290

$ws.Range("H107").Value = 1111
$ws.Range("I107").Value = 1111
$ws.Range("K107").Value = 1111
$ws.Range("M107").Value = 809

$ws.Range("H122").Value = 11942.889
$ws.Range("I122").Value = 7926.5713
$ws.Range("K122").Value = 23779.7139
$ws.Range("M122").Value = -21329.7139

$ws.Range("H132").Value = 5977.8
$ws.Range("I132").Value = 2475.7
$ws.Range("K132").Value = 7427.099999999999
$ws.Range("M132").Value = -4897.099999999999

$ws.Range("H136").Value = 4295.316
$ws.Range("I136").Value = 2704
$ws.Range("J136").Value = 4593.6875
$ws.Range("K136").Value = 8112
$ws.Range("L136").Value = 13781.0625
$ws.Range("M136").Value = -5562
$ws.Range("N136").Value = -18881.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 411.5
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 463.8
$ws.Range("K5").Value = 450
$ws.Range("L5").Value = 1391.4
$ws.Range("N5").Value = -1615.4
$ws.Range("M5").Value = -338

$ws.Range("H6").Value = 48.2
$ws.Range("I6").Value = 42.5
$ws.Range("J6").Value = 56.75
$ws.Range("K6").Value = 127.5
$ws.Range("L6").Value = 170.25
$ws.Range("M6").Value = -14.5
$ws.Range("N6").Value = -396.25

$ws.Range("H23").Value = 393
$ws.Range("J23").Value = 495
$ws.Range("L23").Value = 1485
$ws.Range("N23").Value = -1955

$ws.Range("H122").Value = 549.375
$ws.Range("I122").Value = 319.2
$ws.Range("K122").Value = 2872.8
$ws.Range("M122").Value = -422.7999999999997

$ws.Range("H135").Value = 411.5
$ws.Range("I135").Value = 150
$ws.Range("J135").Value = 463.8
$ws.Range("K135").Value = 1350
$ws.Range("L135").Value = 4174.2
$ws.Range("N135").Value = -9244.200000000001
$ws.Range("M135").Value = 1185

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 575.3
$ws.Range("I2").Value = 1109
$ws.Range("K2").Value = 1109
$ws.Range("M2").Value = -996

$ws.Range("H98").Value = 26876.572
$ws.Range("J98").Value = 26876.572
$ws.Range("L98").Value = 26876.572
$ws.Range("N98").Value = -32866.572

$ws.Range("H113").Value = 4539.8887
$ws.Range("I113").Value = 4539.8887
$ws.Range("K113").Value = 4539.8887
$ws.Range("M113").Value = -2369.8887

$ws.Range("H126").Value = 2900
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -13640

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1000
$ws.Range("J7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("N7").Value = -1224

$ws.Range("H16").Value = 1665.8462
$ws.Range("I16").Value = 1667.7
$ws.Range("K16").Value = 1667.7
$ws.Range("M16").Value = -1497.7

$ws.Range("H22").Value = 1450
$ws.Range("I22").Value = 2000
$ws.Range("K22").Value = 2000
$ws.Range("M22").Value = -1705

$ws.Range("H27").Value = 1450
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1893

$ws.Range("H35").Value = 3001
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = ""

$ws.Range("H40").Value = 8600.333000000001
$ws.Range("I40").Value = 8600.333000000001
$ws.Range("K40").Value = 8600.333000000001
$ws.Range("M40").Value = -8464.333000000001

$ws.Range("H58").Value = 12300
$ws.Range("I58").Value = 13450
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 13450
$ws.Range("L58").Value = 10000
$ws.Range("N58").Value = -10520
$ws.Range("M58").Value = -13190

$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -798

$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170

$ws.Range("H126").Value = 1000
$ws.Range("J126").Value = 1000
$ws.Range("L126").Value = 3000
$ws.Range("N126").Value = -7940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180

$ws.Range("H132").Value = 6611.3335
$ws.Range("I132").Value = 6933.8
$ws.Range("K132").Value = 20801.4
$ws.Range("M132").Value = -18271.4
